$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @(0.0132, 0.0396, 0.0396, 0.05280000000000001, 0.066)
  3  = @(0.0132, 0.0198, 0.0132, 0.0264, 0.0462)
  4  = @(0.0066, 0.0132, 0.0132, 0.033, 0.0264)
  5  = @(0.0198, 0, 0.0066, 0.0198, 0.0264)
  6  = @(0.0066, 0.0132, 0.0132, 0, 0.033)
  7  = @(0.0198, 0.0396, 0.05280000000000001, 0.0198, 0.09899999999999998)
  8  = @(0.0198, 0.0264, 0.0396, 0.0198, 0.08579999999999999)
  9  = @(0.0264, 0.0198, 0.033, 0.0462, 0.0726)
  10 = @(0.0462, 0.0264, 0.0396, 0.0462, 0.1188)
  11 = @(0.0198, 0.0396, 0.0396, 0.07919999999999999, 0.09239999999999998)
  12 = @(0.0264, 0.1188, 0.132, 0.1188, 0.1385999999999999)
  13 = @(0.0726, 0.1122, 0.08579999999999999, 0.09239999999999998, 0.2837999999999998)
  14 = @(0.07919999999999999, 0.132, 0.08579999999999999, 0.0726, 0.2309999999999999)
  15 = @(0.0462, 0.1188, 0.132, 0.09899999999999998, 0.2705999999999998)
  16 = @(0.066, 0.1056, 0.1254, 0.1188, 0.2507999999999999)
  17 = @(0.07919999999999999, 0.09239999999999998, 0.1913999999999999, 0.1517999999999999, 0.2177999999999999)
  18 = @(0.0726, 0.1188, 0.1781999999999999, 0.1385999999999999, 0.2903999999999998)
  19 = @(0.1056, 0.1649999999999999, 0.2045999999999999, 0.1847999999999999, 0.4025999999999997)
  20 = @(0.0726, 0.1451999999999999, 0.1517999999999999, 0.1254, 0.3365999999999998)
  21 = @(0.066, 0.1847999999999999, 0.2243999999999999, 0.1517999999999999, 0.4157999999999997)
  22 = @(0.05940000000000001, 0.1517999999999999, 0.1583999999999999, 0.07919999999999999, 0.3497999999999998)
  23 = @(0.0726, 0.1056, 0.132, 0.09899999999999998, 0.2639999999999998)
  24 = @(0.0462, 0.05940000000000001, 0.0726, 0.0726, 0.1451999999999999)
  25 = @(0.0264, 0.0396, 0.0462, 0.0264, 0.09239999999999998)
}

$cols = @("T", "U", "V", "W", "X")

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])$row").Value = $vals[$i]
  }
}
